# Weekly update: insert a new Pimiento price record as row 46, pushing the
# existing rows 46-71 down to 47-72 (dimension grows from R71 to R72).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 46 (shifts rows 46..71 -> 47..72).
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with the latest observation.
$ws.Cells.Item(46, 1).Value2  = 12
$ws.Cells.Item(46, 2).Value2  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(46, 3).Value2  = "Metropolitana"
$ws.Cells.Item(46, 4).Value2  = 44460
$ws.Cells.Item(46, 5).Value2  = 13
$ws.Cells.Item(46, 6).Value2  = 100112002
$ws.Cells.Item(46, 7).Value2  = "Pimiento"
$ws.Cells.Item(46, 8).Value2  = "Zafiro verde"
$ws.Cells.Item(46, 9).Value2  = "Primera"
$ws.Cells.Item(46, 10).Value2 = 55
$ws.Cells.Item(46, 11).Value2 = 36000
$ws.Cells.Item(46, 12).Value2 = 36000
$ws.Cells.Item(46, 13).Value2 = 36000
$ws.Cells.Item(46, 14).Value2 = "$/caja 18 kilos"
$ws.Cells.Item(46, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(46, 16).Value2 = 2000
$ws.Cells.Item(46, 17).Value2 = 18
$ws.Cells.Item(46, 18).Value2 = "Hortaliza"
